# Form the consolidated report: populate the "Absent" (column H) values
# for each attendance row that were left blank/zero previously.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H16").Value = 0
